$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(97,1).Value = 44839.02611079861
$ws.Cells.Item(97,2).Value = "KAPALI"
$ws.Cells.Item(97,3).Value = "Sistem Kapali"

Write-Host "Done"
